# Add the new "2021-10-29" monthly results sheet as the first sheet in the
# workbook. The new sheet carries the same shape/content as the most recent
# existing monthly sheet ("2021-10-26"), so we duplicate that sheet and move
# the duplicate to the front, then rename it. This leaves every other
# existing sheet's name, order (shifted by one), and data untouched.

$wb = $excel.ActiveWorkbook

$source = $wb.Worksheets.Item("2021-10-26")
$firstSheet = $wb.Worksheets.Item(1)

# Worksheet.Copy(Before) inserts the copy immediately before "Before" and
# shifts every other sheet one position to the right.
$source.Copy($firstSheet)

$newSheet = $wb.Worksheets.Item(1)
$newSheet.Name = "2021-10-29"
